$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E values regenerated to filter save games (G = sum of B:E)
$data = @{
    2 = @(0.06328177979961902, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569)
    3 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
    5 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569)
    6 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $sum = $vals[0] + $vals[1] + $vals[2] + $vals[3]
    $ws.Cells.Item($row, 7).Value = $sum
}
